# completed March 2020 TP2 run 11 and 7 samples of run 12
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Run 11's last sample (row 48) actually belongs with the CRM batch that was
# opened 20210721 (batch value 2224.47, not 2225.47 as originally entered).
$ws.Range("C48").Value = 2224.4699999999998
$ws.Range("D48").Formula = "=100*(B48-C48)/C48"

# New sample from run 12 (row 49), same CRM batch.
$ws.Range("A49").Value = 20211110
$ws.Range("B49").Value = 2220.49470192943
$ws.Range("C49").Value = 2224.4699999999998
$ws.Range("D49").Formula = "=100*(B49-C49)/C49"
$ws.Range("E49").Value = 181
$ws.Range("F49").Value = "CRM OPENED 20210721"

# Leave the selection where the author left it after entering the new row.
$ws.Range("G48").Select()
